# Apply fix to extract.py output: correct headers (drop Age Group columns)
# and trailing columns (drop Certifications/Alignment/Assurance), and replace
# sample/placeholder rows with corrected Singtel data (rows 2-13).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove "Current Employees by Age Groups (Millennials %)" and
# "New Hires and Turnover by Age Groups (Millennials %)" columns (L:M)
$ws.Columns("L:M").Delete()

# Remove trailing "List of Relevant Certifications", "Alignment with Frameworks
# and Disclosure Practices" and "Assurance of Sustainability Report" columns (X:AB)
$ws.Columns("X:AB").Delete()

# Row 2
$ws.Range("A2").Value = "Singtel"
$ws.Range("B2").Value = 2012
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2:W2").ClearContents()

# Row 3
$ws.Range("A3").Value = "Singtel"
$ws.Range("B3").Value = 2013
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("F3:W3").ClearContents()

# Row 4
$ws.Range("A4").Value = "Singtel"
$ws.Range("B4").Value = 2015
$ws.Range("C4:W4").ClearContents()

# Row 5
$ws.Range("A5").Value = "Singtel"
$ws.Range("B5").Value = 2016
$ws.Range("C5:W5").ClearContents()

# Row 6
$ws.Range("A6").Value = "Singtel"
$ws.Range("B6").Value = 2017
$ws.Range("C6:W6").ClearContents()

# Row 7
$ws.Range("A7").Value = "Singtel"
$ws.Range("B7").Value = 2018
$ws.Range("C7").Value = 4085
$ws.Range("D7").Value = 154152
$ws.Range("E7").Value = 6392
$ws.Range("F7").Value = 164629
$ws.Range("G7").Value = 374193.077132
$ws.Range("H7").Value = 620.864
$ws.Range("I7").Value = 7538
$ws.Range("J7:W7").ClearContents()

# Row 8
$ws.Range("A8").Value = "Singtel"
$ws.Range("B8").Value = 2019
$ws.Range("C8").Value = 3741
$ws.Range("D8").Value = 153650
$ws.Range("E8").Value = 5175
$ws.Range("F8").Value = 164629
$ws.Range("G8").Value = 1347094
$ws.Range("H8").Value = 683.847
$ws.Range("I8").Value = 7658
$ws.Range("J8").ClearContents()
$ws.Range("K8").Value = 5
$ws.Range("L8").Value = 18.3
$ws.Range("M8").Value = 12589
$ws.Range("N8").ClearContents()
$ws.Range("O8").Value = 0
$ws.Range("P8:W8").ClearContents()

# Row 9
$ws.Range("A9").Value = "Singtel"
$ws.Range("B9").Value = 2020
$ws.Range("C9").Value = 5749
$ws.Range("D9").Value = 158687
$ws.Range("E9").Value = 5500
$ws.Range("F9").Value = 165331
$ws.Range("G9").Value = 1466802
$ws.Range("H9").Value = 683.847
$ws.Range("I9").Value = 7658
$ws.Range("J9").Value = 45
$ws.Range("K9").Value = 782
$ws.Range("L9").Value = 15.7
$ws.Range("M9").Value = 22914
$ws.Range("N9").Value = 40.2
$ws.Range("O9").Value = 0
$ws.Range("P9").ClearContents()
$ws.Range("Q9").Value = 6
$ws.Range("R9").Value = 5
$ws.Range("S9").ClearContents()
$ws.Range("T9").Value = 30
$ws.Range("U9").Value = 25
$ws.Range("V9").ClearContents()
$ws.Range("W9").Value = 1000

# Row 10
$ws.Range("A10").Value = "Singtel"
$ws.Range("B10").Value = 2021
$ws.Range("C10").Value = 9500
$ws.Range("D10").Value = 18000
$ws.Range("E10").Value = 22000
$ws.Range("F10").Value = 165331
$ws.Range("G10").Value = 1602698
$ws.Range("H10").Value = 3000
$ws.Range("I10").Value = 4150
$ws.Range("J10").Value = 45
$ws.Range("K10").Value = 439
$ws.Range("L10").Value = 13.2
$ws.Range("M10").Value = 12391
$ws.Range("N10").Value = 48.3
$ws.Range("O10").Value = 0
$ws.Range("P10").ClearContents()
$ws.Range("Q10").Value = 15
$ws.Range("R10").Value = 2
$ws.Range("S10").ClearContents()
$ws.Range("T10").Value = 30
$ws.Range("U10").Value = 25
$ws.Range("V10").ClearContents()
$ws.Range("W10").Value = 95

# Row 11
$ws.Range("A11").Value = "Singtel"
$ws.Range("B11").Value = 2022
$ws.Range("C11").Value = 5000
$ws.Range("D11").Value = 7000
$ws.Range("E11").Value = 3000
$ws.Range("F11").Value = 15000
$ws.Range("G11").Value = 120000
$ws.Range("H11").Value = 500
$ws.Range("I11").Value = 1200
$ws.Range("J11").Value = 40
$ws.Range("K11").ClearContents()
$ws.Range("L11").Value = 10
$ws.Range("M11").ClearContents()
$ws.Range("N11").Value = 20
$ws.Range("O11").Value = 0
$ws.Range("P11").ClearContents()
$ws.Range("Q11").Value = 11
$ws.Range("R11").ClearContents()
$ws.Range("S11").Value = 50
$ws.Range("T11:W11").ClearContents()

# Row 12
$ws.Range("A12").Value = "Singtel"
$ws.Range("B12").Value = 2030
$ws.Range("C12:W12").ClearContents()

# Row 13
$ws.Range("A13").Value = "Singtel"
$ws.Range("B13").Value = 2050
$ws.Range("C13:W13").ClearContents()

Write-Host "Edit applied"